$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 4.014144333333333
$ws.Cells.Item(2, 8).Value = 12.042433
$ws.Cells.Item(2, 9).Value = 0.0496505710047397
$ws.Cells.Item(2, 10).Value = 0.0496505710047397
$ws.Cells.Item(2, 13).Value = 3.303267
$ws.Cells.Item(2, 14).Value = 9.909801000000002
$ws.Cells.Item(2, 15).Value = 0.03362563178859915
$ws.Cells.Item(2, 16).Value = 0.03362563178859915
$ws.Cells.Item(2, 17).Value = 13.259790509537
$ws.Cells.Item(2, 18).Value = 119.338114585833
$ws.Cells.Item(2, 19).Value = 0.001669531818699075
$ws.Cells.Item(2, 20).Value = 0.001669531818699075

$ws.Cells.Item(3, 7).Value = 4.014144333333333
$ws.Cells.Item(3, 8).Value = 12.042433
$ws.Cells.Item(3, 9).Value = 0.0496505710047397
$ws.Cells.Item(3, 10).Value = 0.0496505710047397
$ws.Cells.Item(3, 13).Value = 37.82684066666667
$ws.Cells.Item(3, 15).Value = 0.3850586149964086
$ws.Cells.Item(3, 16).Value = 0.3850586149964086
$ws.Cells.Item(3, 17).Value = 151.8423981100029
$ws.Cells.Item(3, 18).Value = 1366.581582990026
$ws.Cells.Item(3, 19).Value = 0.01911838010486592
$ws.Cells.Item(3, 20).Value = 0.01911838010486592

$ws.Cells.Item(4, 7).Value = 4.014144333333333
$ws.Cells.Item(4, 8).Value = 12.042433
$ws.Cells.Item(4, 9).Value = 0.0496505710047397
$ws.Cells.Item(4, 10).Value = 0.0496505710047397
$ws.Cells.Item(4, 13).Value = 9.149396
$ws.Cells.Item(4, 14).Value = 27.448188
$ws.Cells.Item(4, 15).Value = 0.09313634682999644
$ws.Cells.Item(4, 16).Value = 0.09313634682999644
$ws.Cells.Item(4, 17).Value = 36.72699610682266
$ws.Cells.Item(4, 18).Value = 330.5429649614039
$ws.Cells.Item(4, 19).Value = 0.004624272801404801
$ws.Cells.Item(4, 20).Value = 0.004624272801404801

$ws.Cells.Item(5, 7).Value = 4.014144333333333
$ws.Cells.Item(5, 8).Value = 12.042433
$ws.Cells.Item(5, 9).Value = 0.0496505710047397
$ws.Cells.Item(5, 10).Value = 0.0496505710047397
$ws.Cells.Item(5, 13).Value = 47.95707433333333
$ws.Cells.Item(5, 14).Value = 143.871223
$ws.Cells.Item(5, 15).Value = 0.4881794063849957
$ws.Cells.Item(5, 16).Value = 0.4881794063849957
$ws.Cells.Item(5, 17).Value = 192.5066181783954
$ws.Cells.Item(5, 18).Value = 1732.559563605559
$ws.Cells.Item(5, 19).Value = 0.02423838627976991
$ws.Cells.Item(5, 20).Value = 0.02423838627976991

$ws.Cells.Item(6, 9).Value = 0.6580818039484316
$ws.Cells.Item(6, 10).Value = 0.6580818039484316
$ws.Cells.Item(6, 13).Value = 3.303267
$ws.Cells.Item(6, 14).Value = 9.909801000000002
$ws.Cells.Item(6, 15).Value = 0.03362563178859915
$ws.Cells.Item(6, 16).Value = 0.03362563178859915
$ws.Cells.Item(6, 17).Value = 175.748771502777
$ws.Cells.Item(6, 18).Value = 1581.738943524993
$ws.Cells.Item(6, 19).Value = 0.02212841642634706
$ws.Cells.Item(6, 20).Value = 0.02212841642634706

$ws.Cells.Item(7, 9).Value = 0.6580818039484316
$ws.Cells.Item(7, 10).Value = 0.6580818039484316
$ws.Cells.Item(7, 13).Value = 37.82684066666667
$ws.Cells.Item(7, 15).Value = 0.3850586149964086
$ws.Cells.Item(7, 16).Value = 0.3850586149964086
$ws.Cells.Item(7, 18).Value = 18113.03385193555
$ws.Cells.Item(7, 19).Value = 0.2534000679827212
$ws.Cells.Item(7, 20).Value = 0.2534000679827212

$ws.Cells.Item(8, 9).Value = 0.6580818039484316
$ws.Cells.Item(8, 10).Value = 0.6580818039484316
$ws.Cells.Item(8, 13).Value = 9.149396
$ws.Cells.Item(8, 14).Value = 27.448188
$ws.Cells.Item(8, 15).Value = 0.09313634682999644
$ws.Cells.Item(8, 16).Value = 0.09313634682999644
$ws.Cells.Item(8, 17).Value = 486.7893231132759
$ws.Cells.Item(8, 18).Value = 4381.103908019483
$ws.Cells.Item(8, 19).Value = 0.06129133513505085
$ws.Cells.Item(8, 20).Value = 0.06129133513505085

$ws.Cells.Item(9, 9).Value = 0.6580818039484316
$ws.Cells.Item(9, 10).Value = 0.6580818039484316
$ws.Cells.Item(9, 13).Value = 47.95707433333333
$ws.Cells.Item(9, 14).Value = 143.871223
$ws.Cells.Item(9, 15).Value = 0.4881794063849957
$ws.Cells.Item(9, 16).Value = 0.4881794063849957
$ws.Cells.Item(9, 17).Value = 2551.533648037137
$ws.Cells.Item(9, 18).Value = 22963.80283233423
$ws.Cells.Item(9, 19).Value = 0.3212619844043125
$ws.Cells.Item(9, 20).Value = 0.3212619844043125

$ws.Cells.Item(10, 7).Value = 20.763658
$ws.Cells.Item(10, 8).Value = 62.29097400000001
$ws.Cells.Item(10, 9).Value = 0.256823718889812
$ws.Cells.Item(10, 10).Value = 0.256823718889812
$ws.Cells.Item(10, 13).Value = 3.303267
$ws.Cells.Item(10, 14).Value = 9.909801000000002
$ws.Cells.Item(10, 15).Value = 0.03362563178859915
$ws.Cells.Item(10, 16).Value = 0.03362563178859915
$ws.Cells.Item(10, 17).Value = 68.58790627068602
$ws.Cells.Item(10, 18).Value = 617.2911564361741
$ws.Cells.Item(10, 19).Value = 0.008635859805967514
$ws.Cells.Item(10, 20).Value = 0.008635859805967514

$ws.Cells.Item(11, 7).Value = 20.763658
$ws.Cells.Item(11, 8).Value = 62.29097400000001
$ws.Cells.Item(11, 9).Value = 0.256823718889812
$ws.Cells.Item(11, 10).Value = 0.256823718889812
$ws.Cells.Item(11, 13).Value = 37.82684066666667
$ws.Cells.Item(11, 15).Value = 0.3850586149964086
$ws.Cells.Item(11, 16).Value = 0.3850586149964086
$ws.Cells.Item(11, 17).Value = 785.4235828231589
$ws.Cells.Item(11, 18).Value = 7068.812245408429
$ws.Cells.Item(11, 19).Value = 0.098892185493938
$ws.Cells.Item(11, 20).Value = 0.098892185493938

$ws.Cells.Item(12, 7).Value = 20.763658
$ws.Cells.Item(12, 8).Value = 62.29097400000001
$ws.Cells.Item(12, 9).Value = 0.256823718889812
$ws.Cells.Item(12, 10).Value = 0.256823718889812
$ws.Cells.Item(12, 13).Value = 9.149396
$ws.Cells.Item(12, 14).Value = 27.448188
$ws.Cells.Item(12, 15).Value = 0.09313634682999644
$ws.Cells.Item(12, 16).Value = 0.09313634682999644
$ws.Cells.Item(12, 17).Value = 189.974929450568
$ws.Cells.Item(12, 18).Value = 1709.774365055112
$ws.Cells.Item(12, 19).Value = 0.02391962295669104
$ws.Cells.Item(12, 20).Value = 0.02391962295669104

$ws.Cells.Item(13, 7).Value = 20.763658
$ws.Cells.Item(13, 8).Value = 62.29097400000001
$ws.Cells.Item(13, 9).Value = 0.256823718889812
$ws.Cells.Item(13, 10).Value = 0.256823718889812
$ws.Cells.Item(13, 13).Value = 47.95707433333333
$ws.Cells.Item(13, 14).Value = 143.871223
$ws.Cells.Item(13, 15).Value = 0.4881794063849957
$ws.Cells.Item(13, 16).Value = 0.4881794063849957
$ws.Cells.Item(13, 17).Value = 995.7642901379114
$ws.Cells.Item(13, 18).Value = 8961.878611241202
$ws.Cells.Item(13, 19).Value = 0.1253760506332154
$ws.Cells.Item(13, 20).Value = 0.1253760506332154

$ws.Cells.Item(14, 7).Value = 2.865565333333334
$ws.Cells.Item(14, 8).Value = 8.596696000000001
$ws.Cells.Item(14, 9).Value = 0.03544390615701676
$ws.Cells.Item(14, 10).Value = 0.03544390615701676
$ws.Cells.Item(14, 13).Value = 3.303267
$ws.Cells.Item(14, 14).Value = 9.909801000000002
$ws.Cells.Item(14, 15).Value = 0.03362563178859915
$ws.Cells.Item(14, 16).Value = 0.03362563178859915
$ws.Cells.Item(14, 17).Value = 9.465727401944003
$ws.Cells.Item(14, 18).Value = 85.19154661749603
$ws.Cells.Item(14, 19).Value = 0.001191823737585508
$ws.Cells.Item(14, 20).Value = 0.001191823737585508

$ws.Cells.Item(15, 7).Value = 2.865565333333334
$ws.Cells.Item(15, 8).Value = 8.596696000000001
$ws.Cells.Item(15, 9).Value = 0.03544390615701676
$ws.Cells.Item(15, 10).Value = 0.03544390615701676
$ws.Cells.Item(15, 13).Value = 37.82684066666667
$ws.Cells.Item(15, 15).Value = 0.3850586149964086
$ws.Cells.Item(15, 16).Value = 0.3850586149964086
$ws.Cells.Item(15, 17).Value = 108.3952832839236
$ws.Cells.Item(15, 18).Value = 975.5575495553122
$ws.Cells.Item(15, 19).Value = 0.01364798141488356
$ws.Cells.Item(15, 20).Value = 0.01364798141488356

$ws.Cells.Item(16, 7).Value = 2.865565333333334
$ws.Cells.Item(16, 8).Value = 8.596696000000001
$ws.Cells.Item(16, 9).Value = 0.03544390615701676
$ws.Cells.Item(16, 10).Value = 0.03544390615701676
$ws.Cells.Item(16, 13).Value = 9.149396
$ws.Cells.Item(16, 14).Value = 27.448188
$ws.Cells.Item(16, 15).Value = 0.09313634682999644
$ws.Cells.Item(16, 16).Value = 0.09313634682999644
$ws.Cells.Item(16, 17).Value = 26.21819199853867
$ws.Cells.Item(16, 18).Value = 235.963727986848
$ws.Cells.Item(16, 19).Value = 0.00330111593684976
$ws.Cells.Item(16, 20).Value = 0.00330111593684976

$ws.Cells.Item(17, 7).Value = 2.865565333333334
$ws.Cells.Item(17, 8).Value = 8.596696000000001
$ws.Cells.Item(17, 9).Value = 0.03544390615701676
$ws.Cells.Item(17, 10).Value = 0.03544390615701676
$ws.Cells.Item(17, 13).Value = 47.95707433333333
$ws.Cells.Item(17, 14).Value = 143.871223
$ws.Cells.Item(17, 15).Value = 0.4881794063849957
$ws.Cells.Item(17, 16).Value = 0.4881794063849957
$ws.Cells.Item(17, 17).Value = 137.4241296976898
$ws.Cells.Item(17, 18).Value = 1236.817167279208
$ws.Cells.Item(17, 19).Value = 0.01730298506769794
$ws.Cells.Item(17, 20).Value = 0.01730298506769794
